# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "239.32"); force
# the whole data range to Text format first so Excel does not silently
# convert these inline strings into real numbers on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '29.506.30'
$ws.Range('E2').Value = '  +0.95%  '
$ws.Range('D3').Value = '1.871.59'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '0.7175'
$ws.Range('E5').Value = '  +1.14%  '
$ws.Range('D6').Value = '239.32'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('D8').Value = '0.07827'
$ws.Range('E8').Value = '  -4.25%  '
$ws.Range('D9').Value = '0.3070'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').Value = '25.32'
$ws.Range('E10').Value = '  +8.80%  '
$ws.Range('D11').Value = '0.08227'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '1.896.39'
$ws.Range('E12').Value = '  +1.50%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '0.7222'
$ws.Range('E13').Value = '  +1.91%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.231'
$ws.Range('E14').Value = '  +1.13%  '
$ws.Range('D15').Value = '89.90'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '29.543.43'
$ws.Range('E16').Value = '  +1.05%  '
$ws.Range('D17').Value = '5.822'
$ws.Range('E17').Value = '  +0.63%  '
$ws.Range('D18').Value = '0.000007847'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = '241.10'
$ws.Range('E19').Value = '  +1.67%  '
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('D21').Value = '2.130.49'
$ws.Range('E21').Value = '  +1.21%  '
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '7.726'
$ws.Range('E24').Value = '  +4.33%  '
$ws.Range('D25').Value = '0.1560'
$ws.Range('E25').Value = '  +7.41%  '
$ws.Range('D26').Value = '162.66'
$ws.Range('D27').Value = '8.942'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('D29').Value = '1.936'
$ws.Range('E29').Value = '  -1.20%  '
$ws.Range('D30').Value = '1.360'
$ws.Range('E30').Value = '  -4.59%  '
$ws.Range('E31').Value = '  -0.26%  '
$ws.Range('E32').Value = '  -1.28%  '
$ws.Range('D33').Value = '4.073'
$ws.Range('E33').Value = '  +1.10%  '
$ws.Range('D34').Value = '0.05250'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +2.29%  '
$ws.Range('D36').Value = '0.7156'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('D37').Value = '1.002'
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('D38').Value = '2.671'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '0.01868'
$ws.Range('E39').Value = '  +0.64%  '
$ws.Range('E40').Value = '  -0.45%  '
$ws.Range('D41').Value = '1.174.66'
$ws.Range('E41').Value = '  +2.54%  '
$ws.Range('D42').Value = '0.9082'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = '5.998'
$ws.Range('E43').Value = '  +1.84%  '
$ws.Range('D44').Value = '0.4306'
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('D45').Value = '71.61'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').Value = '102.52'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('D48').Value = '0.5363'
$ws.Range('E48').Value = '  -0.83%  '
$ws.Range('D49').Value = '1.765'
$ws.Range('E49').Value = '  -0.56%  '
$ws.Range('D50').Value = '9.167'
$ws.Range('E50').Value = '  -0.52%  '
$ws.Range('D51').Value = '7.021'
$ws.Range('E51').Value = '  +0.93%  '
